$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append: row, date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg per 100mila ab.
$data = @(
    @(358, 44432, 1, 8, 70.86544423775356),
    @(359, 44433, 0, 8, 70.86544423775356),
    @(360, 44434, 1, 8, 70.86544423775356),
    @(361, 44435, 3, 10, 88.58180529719195),
    @(362, 44436, 1, 8, 70.86544423775356),
    @(363, 44437, 1, 7, 62.00726370803437),
    @(364, 44438, 2, 9, 79.72362476747276),
    @(365, 44439, 1, 9, 79.72362476747276),
    @(366, 44440, 0, 9, 79.72362476747276)
)

# Use the formatting of the last existing data row (357) as the template for column A's style.
$ws.Cells.Item(357, 1).Copy()

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
